$d = $word.ActiveDocument

# --- Change 1: fix hyphen in "документо-ориентированные" (inside a bookmark) ---
$d.Content.Find.Execute("документо-ориентированные", $true, $false, $false, $false, $false, $true, 1, $false, "документоориентированные", 2)

# --- Change 2: insert two new paragraphs about document-oriented DBs right after
#     the paragraph ending in "[Hoffner]" and before the "К концу 1980-х годов..." paragraph ---
$rng = $d.Content
$rng.Find.Execute("[Hoffner]", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Collapse(0)
$rng.InsertParagraphAfter()
$rng.MoveStart(1, 1)
$rng.Collapse(1)
$rng.InsertAfter("Документоориентированная СУБД хранит данные в виде структурированных документов, обычно в формате XML или JSON. ")
$rng.Collapse(0)
$rng.InsertAfter("При этом определение «документоориентированная СУБД» не подразумевает какую-либо  специфику насчёт модели хранения: документоориентированные СУБД могут выполнять ACID-транзакции или другие функции традиционных реляционных СУБД, хотя популярные документоориентированные обеспечивают относительно скромную транзакционную поддержку.")

$rng.Collapse(0)
$rng.InsertParagraphAfter()
$rng.MoveStart(1, 1)
$rng.Collapse(1)
$rng.InsertAfter("Документоориентированные базы данных, позволяя описывать данные без использования схемы, возможно, являются золотой серединой между жёсткой схемой реляционных баз данных и свободных от схемы хранилищ «ключ-значение». Сочетание с практикой веб-разработки вылилось в появление JSON-баз данных (MongoDB  в частности), которые стали выбором по умолчанию для многих веб-разработчиков ")
$rng.Collapse(0)
$rng.InsertAfter("[Harrison G. Next Gen.].")

# --- Change 2b: prepend a new sentence about the relational approach's origin as its
#     own paragraph, right before the "К концу 1980-х годов..." paragraph ---
$rng.Collapse(0)
$rng.InsertParagraphAfter()
$rng.MoveStart(1, 1)
$rng.Collapse(1)
$rng.InsertAfter("Реляционный подход к СУБД зародился в конце 1960-х годов. ")

# --- Change 3: split "Кузнецов Основы баз данных" into "Кузнецов Основы " + "БД" ---
$d.Content.Find.Execute("Кузнецов Основы баз данных", $true, $false, $false, $false, $false, $true, 1, $false, "Кузнецов Основы БД", 2)
